# Apply the "postman environment" update to sampletest.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "c-demo_ui": add a new config row (row 34) describing the
# postman environment filename used by the test data reader.
# ---------------------------------------------------------------------
$wsDemo = $wb.Worksheets.Item("c-demo_ui")

$wsDemo.Range("A34").Value = "postman.environment.filename"
$wsDemo.Range("B34").Value = "PostmanEnvSample/postman_environment.json"

# Select the new row (matches the saved selection/scroll state).
$wsDemo.Activate()
$wsDemo.Range("B36").Select()

# ---------------------------------------------------------------------
# Sheet "t-omdb&imdb": wire the query-string values used to resolve
# each movie id from the postman environment (apikey + movie id).
# ---------------------------------------------------------------------
$wsOmdb = $wb.Worksheets.Item("t-omdb&imdb")

$wsOmdb.Range("F3").Value = "query::apikey::{{apikey}}"
$wsOmdb.Range("G3").Value = "query::i::{{wall-eId}}"
$wsOmdb.Range("G3").Font.Bold = $false

$wsOmdb.Range("F4").Value = "query::apikey::{{apikey}}"
$wsOmdb.Range("G4").Value = "query::i::{{lionKingId}}"
$wsOmdb.Range("G4").Font.Bold = $false

$wsOmdb.Range("F5").Value = "query::apikey::{{apikey}}"
$wsOmdb.Range("G5").Value = "query::i::{{shaunId}}"
$wsOmdb.Range("G5").Font.Bold = $false

# Widen column G so the new values are fully visible.
$wsOmdb.Columns.Item(7).ColumnWidth = 34.6667

# Update the selection for this sheet (it is the active sheet on save).
$wsOmdb.Activate()
$wsOmdb.Range("G5").Select()
